# Add two new columns, I (I0) and J (IF), to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells -----------------------------------------------------
# Set the header text first.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, borders, centered alignment) from the
# existing "IP" header cell (H1) onto the two new header cells so they
# match the look of the other headers.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# --- Data rows ----------------------------------------------------------
$data = @(
    @(2, 3),
    @(4, 5),
    @(8, 9),
    @(1, 4),
    @(1, 4),
    @(1, 5),
    @(1, 6),
    @(6, 6),
    @(5, 9),
    @(1, 4),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 4),
    @(1, 5),
    @(4, 5),
    @(7, 8),
    @(8, 8),
    @(7, 8),
    @(12, 12),
    @(7, 7),
    @(1, 3),
    @(1, 5),
    @(1, 4),
    @(1, 4),
    @(1, 4),
    @(1, 4),
    @(3, 4),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Range("I$row").Value = $data[$i][0]
    $ws.Range("J$row").Value = $data[$i][1]
}
